# Generate Report for Handback
# The "fed47d42-f980-4d67-81ac-105323179322" file has moved from
# "Ready for handoff" to "Handed back: in sync with en-US", and its
# handback timestamps were refreshed accordingly on the per-language
# report sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: row for fed47d42-... (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet: status + refreshed handback datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $newStatus
$wsZhCn.Range("G2").Value = "2016-03-01 09:52:04"
$wsZhCn.Range("G3").Value = "2016-03-01 09:52:04"

# --- de-de sheet: status + refreshed handback datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $newStatus
$wsDeDe.Range("G2").Value = "2016-03-01 09:52:22"
$wsDeDe.Range("G3").Value = "2016-03-01 09:52:22"
